$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 412, shifting existing rows 412:517 down to 413:518
$ws.Rows.Item(412).Insert()

# Populate the newly inserted row 412 with the new record
$ws.Cells.Item(412, 1).Value = 8
$ws.Cells.Item(412, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(412, 3).Value = "Coquimbo"
$ws.Cells.Item(412, 4).Value = 44855
$ws.Cells.Item(412, 5).Value = 4
$ws.Cells.Item(412, 6).Value = 100114001
$ws.Cells.Item(412, 7).Value = "Papa"
$ws.Cells.Item(412, 8).Value = "Cardinal"
$ws.Cells.Item(412, 9).Value = "1a (cosecha)"
$ws.Cells.Item(412, 10).Value = 2000
$ws.Cells.Item(412, 11).Value = 11000
$ws.Cells.Item(412, 12).Value = 12000
$ws.Cells.Item(412, 13).Value = 11500
$ws.Cells.Item(412, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(412, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(412, 16).Value = 460
$ws.Cells.Item(412, 17).Value = 25
$ws.Cells.Item(412, 18).Value = "Hortaliza"
